# Apply updated cryptocurrency price/volume figures (refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.887.91'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '3.857.71'
$ws.Range('E3').Value = '  +2.59%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '601.13'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('E6').Value = '  -3.03%  '
$ws.Range('D7').Value = '3.855.67'
$ws.Range('E7').Value = '  +2.61%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('E12').Value = '  -0.23%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '36.93'
$c.ClearFormats()
$ws.Range('E13').Value = '  -3.07%  '
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').Value = '4.503.07'
$ws.Range('D16').Value = '3.848.94'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '69.060.86'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('E18').Value = '  +2.23%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.47'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('E20').Value = '  -0.28%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '17.07'
$c.ClearFormats()
$ws.Range('E21').Value = '  -1.22%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '484.52'
$c.ClearFormats()
$ws.Range('E22').Value = '  -1.92%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.719'
$c.ClearFormats()
$ws.Range('E23').Value = '  -1.45%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.0000161'
$c.ClearFormats()
$ws.Range('E24').Value = '  +5.50%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '83.96'
$c.ClearFormats()
$ws.Range('E25').Value = '  -1.20%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.24'
$c.ClearFormats()
$ws.Range('E26').Value = '  -2.96%  '
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('E28').Value = '  -0.08%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.96'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.53%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.ClearFormats()
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('D32').Value = '4.007.52'
$ws.Range('E32').Value = '  +2.56%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '32.21'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.78%  '
$ws.Range('E34').Value = '  -4.51%  '
$ws.Range('D35').Value = '3.805.42'
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('E38').Value = '  +1.90%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '5.88'
$c.ClearFormats()
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -2.53%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '437.40'
$c.ClearFormats()
$ws.Range('E42').Value = '  +1.71%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.96'
$c.ClearFormats()
$ws.Range('E43').Value = '  -2.27%  '
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('E46').Value = '  +0.01%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '8.39'
$c.ClearFormats()
$ws.Range('E47').Value = '  -1.07%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '143.45'
$c.ClearFormats()
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').Value = '2.840.05'
$ws.Range('E49').Value = '  +1.51%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0359'
$c.ClearFormats()
$ws.Range('E50').Value = '  +1.45%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '25.97'
$c.ClearFormats()
$ws.Range('E51').Value = '  +12.81%  '
